# Weekly update: a new Cilantro price-report row (2022-03-28, serial 44648)
# is inserted at row 17, pushing all subsequent rows (old 17..72) down by
# one (new 18..73). The sheet's used range grows from A1:R72 to A1:R73.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row before the current row 17 — this shifts rows
# 17..72 down to 18..73 (and carries the date-format style along).
$ws.Rows.Item(17).Insert()

# Populate the newly inserted row 17 with the new weekly data point.
$ws.Range("A17").Value = 1
$ws.Range("B17").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C17").Value = "Arica y Parinacota"
$ws.Range("D17").Value = 44648
$ws.Range("E17").Value = 15
$ws.Range("F17").Value = 100112040
$ws.Range("G17").Value = "Cilantro"
$ws.Range("H17").Value = "Sin especificar"
$ws.Range("I17").Value = "Primera"
$ws.Range("J17").Value = 300
$ws.Range("K17").Value = 1800
$ws.Range("L17").Value = 2000
$ws.Range("M17").Value = 1900
$ws.Range("N17").Value = "`$/atado 1,5 a 2 kilos"
$ws.Range("O17").Value = "Región de Arica y Parinacota"
$ws.Range("P17").Value = 950
$ws.Range("Q17").Value = 2
$ws.Range("R17").Value = "Hortaliza"
